$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 1: add header values for new columns P and Q ---
$ws.Range("P1").Value = 14
$ws.Range("Q1").Value = 15
# Reuse the existing header style (bold font, thin border, centered/top) from O1
$ws.Range("O1").Copy()
$ws.Range("P1:Q1").PasteSpecial(-4122) # xlPasteFormats
$excel.CutCopyMode = $false

# --- Rows 2-25: swap I/K and M/O columns, add P and Q columns ---
for ($r = 2; $r -le 25; $r++) {
    $ws.Cells.Item($r, 9).Value  = 2   # I -> 2
    $ws.Cells.Item($r, 11).Value = 1   # K -> 1
    $ws.Cells.Item($r, 13).Value = 2   # M -> 2
    $ws.Cells.Item($r, 15).Value = 1   # O -> 1
    $ws.Cells.Item($r, 16).Value = 2   # P = 2
    $ws.Cells.Item($r, 17).Value = 2   # Q = 2
}
